# Applies the "Fixes during Regression Testing" regression-run update:
# - refreshes the per-sheet run Notes/Date timestamps left behind by the
#   Katalon regression suite
# - updates the AddDeleteRole RoleName test value from "Space Role" to
#   "Prod Role"
# - leaves the workbook with AddDeleteRole as the active/selected sheet
#   (cell H2 selected), matching the state Excel was saved in

$wb = $excel.ActiveWorkbook

# --- CreateUser sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("CreateUser")
$ws.Range("B2").Value = "Thu Nov 20 15:07:04 IST 2025"
$ws.Range("D2").Value = "Thu Nov 20 19:01:35 IST 2025"

# --- FindUser sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("FindUser")
$ws.Range("B2").Value = "Thu Nov 20 15:07:46 IST 2025"
$ws.Range("D2").Value = "Thu Nov 20 19:02:22 IST 2025"

# --- ModifyUser sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("ModifyUser")
$ws.Range("B2").Value = "Thu Nov 20 21:26:52 IST 2025"
$ws.Range("D2").Value = "Thu Nov 20 19:02:53 IST 2025"

# --- ModifyUserPwd sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("ModifyUserPwd")
$ws.Range("B2").Value = "Thu Nov 20 15:09:20 IST 2025"
$ws.Range("D2").Value = "Thu Nov 20 19:03:49 IST 2025"

# --- FindCaseUser sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("FindCaseUser")
$ws.Range("B2").Value = "Thu Nov 20 15:10:39 IST 2025"
$ws.Range("D2").Value = "Thu Nov 20 19:05:16 IST 2025"

# --- AddDeleteRole sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("AddDeleteRole")
$ws.Range("B2").Value = "Thu Nov 20 21:56:43 IST 2025"
$ws.Range("D2").Value = "Thu Nov 20 19:00:09 IST 2025"
$ws.Range("H2").Value = "Prod Role"

# --- SearchRole sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("SearchRole")
$ws.Range("B2").Value = "Thu Nov 20 15:06:26 IST 2025"
$ws.Range("D2").Value = "Thu Nov 20 19:00:56 IST 2025"

# --- Make AddDeleteRole the active sheet with H2 selected ------------------
$ws = $wb.Worksheets.Item("AddDeleteRole")
$ws.Select()
$ws.Range("H2").Select()
